$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "41.514.17"
Set-TextValue $ws.Range("E2") "  -0.52%  "
Set-TextValue $ws.Range("D3") "2.163.67"
Set-TextValue $ws.Range("E3") "  -2.46%  "
Set-TextValue $ws.Range("E4") "  +0.11%  "
Set-TextValue $ws.Range("D5") "237.47"
Set-TextValue $ws.Range("E5") "  -1.70%  "
Set-TextValue $ws.Range("D6") "0.606"
Set-TextValue $ws.Range("E6") "  -3.20%  "
Set-TextValue $ws.Range("D7") "71.26"
Set-TextValue $ws.Range("E7") "  -1.89%  "
Set-TextValue $ws.Range("E8") "  +0.00%  "
Set-TextValue $ws.Range("D9") "0.574"
Set-TextValue $ws.Range("E9") "  -3.61%  "
Set-TextValue $ws.Range("D10") "39.64"
Set-TextValue $ws.Range("E10") "  -6.00%  "
Set-TextValue $ws.Range("D11") "0.0901"
Set-TextValue $ws.Range("E11") "  -5.18%  "
Set-TextValue $ws.Range("D12") "53.98"
Set-TextValue $ws.Range("E12") "  -4.57%  "
Set-TextValue $ws.Range("D13") "0.0999"
Set-TextValue $ws.Range("E13") "  -3.63%  "
Set-TextValue $ws.Range("D14") "6.67"
Set-TextValue $ws.Range("E14") "  -4.00%  "
Set-TextValue $ws.Range("D15") "2.485.56"
Set-TextValue $ws.Range("E15") "  -2.50%  "
Set-TextValue $ws.Range("E16") "  -0.85%  "
Set-TextValue $ws.Range("D17") "2.157.10"
Set-TextValue $ws.Range("E17") "  -2.69%  "
Set-TextValue $ws.Range("D18") "0.780"
Set-TextValue $ws.Range("E18") "  -6.70%  "
Set-TextValue $ws.Range("D19") "41.361.56"
Set-TextValue $ws.Range("E19") "  -0.72%  "
Set-TextValue $ws.Range("E20") "  -4.44%  "
Set-TextValue $ws.Range("D21") "69.56"
Set-TextValue $ws.Range("E21") "  -4.07%  "
Set-TextValue $ws.Range("D22") "5.75"
Set-TextValue $ws.Range("D23") "9.96"
Set-TextValue $ws.Range("E23") "  -9.46%  "
Set-TextValue $ws.Range("D24") "227.39"
Set-TextValue $ws.Range("E24") "  -0.92%  "
Set-TextValue $ws.Range("D25") "1.98"
Set-TextValue $ws.Range("E25") "  -4.20%  "
Set-TextValue $ws.Range("E26") "  -0.12%  "
Set-TextValue $ws.Range("D27") "10.68"
Set-TextValue $ws.Range("E27") "  -6.36%  "
Set-TextValue $ws.Range("D28") "3.26"
Set-TextValue $ws.Range("E28") "  -10.11%  "
Set-TextValue $ws.Range("E29") "  -4.57%  "
Set-TextValue $ws.Range("E30") "  -0.89%  "
Set-TextValue $ws.Range("D31") "171.57"
Set-TextValue $ws.Range("E31") "  +2.58%  "
Set-TextValue $ws.Range("D32") "19.74"
Set-TextValue $ws.Range("E32") "  -3.56%  "
Set-TextValue $ws.Range("D33") "33.24"
Set-TextValue $ws.Range("E33") "  +10.42%  "
Set-TextValue $ws.Range("E34") "  -3.29%  "
Set-TextValue $ws.Range("D35") "5.11"
Set-TextValue $ws.Range("E35") "  -7.48%  "
Set-TextValue $ws.Range("E36") "  -4.13%  "
Set-TextValue $ws.Range("E37") "  -2.78%  "
Set-TextValue $ws.Range("D38") "4.21"
Set-TextValue $ws.Range("E38") "  -1.61%  "
Set-TextValue $ws.Range("D39") "0.0300"
Set-TextValue $ws.Range("E39") "  -0.86%  "
Set-TextValue $ws.Range("D40") "11.99"
Set-TextValue $ws.Range("E40") "  -10.02%  "
Set-TextValue $ws.Range("E41") "  -2.72%  "
Set-TextValue $ws.Range("D42") "5.33"
Set-TextValue $ws.Range("D43") "58.67"
Set-TextValue $ws.Range("E43") "  -8.80%  "
Set-TextValue $ws.Range("E44") "  -4.48%  "
Set-TextValue $ws.Range("D45") "8.35"
Set-TextValue $ws.Range("E45") "  -4.15%  "
Set-TextValue $ws.Range("D46") "0.0954"
Set-TextValue $ws.Range("E46") "  -4.68%  "
Set-TextValue $ws.Range("D47") "96.12"
Set-TextValue $ws.Range("E47") "  -6.75%  "
Set-TextValue $ws.Range("E48") "  -3.01%  "
Set-TextValue $ws.Range("E49") "  -4.74%  "

# Row 50/51 swap: HuobiToken <-> NEARProtocol (with updated D/E values)
Set-TextValue $ws.Range("B50") "NEARProtocol"
Set-TextValue $ws.Range("C50") "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D50") "2.15"
Set-TextValue $ws.Range("E50") "  -7.45%  "

Set-TextValue $ws.Range("B51") "HuobiToken"
Set-TextValue $ws.Range("C51") "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue $ws.Range("D51") "2.62"
Set-TextValue $ws.Range("E51") "  -2.50%  "
